$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58; this pushes the old rows 58-62
# (the two empty placeholder rows and the three summary rows) down to 59-63.
$ws.Rows.Item(58).Insert()

# Fill the two new data rows (58 and 59) with the new entries.
$ws.Range("A58").Value = 2014
$ws.Range("B58").Value = 3
$ws.Range("C58").Value = 11
$ws.Range("D58").Value = 0.57291666666666663
$ws.Range("E58").Value = 0.66666666666666663

$ws.Range("A59").Value = 2014
$ws.Range("B59").Value = 3
$ws.Range("C59").Value = 11
$ws.Range("D59").Value = 0.79861111111111116
$ws.Range("E59").Value = 0.91666666666666663

# Assign the formulas as a single range write so the two new rows share
# one formula group (matching how the adjacent F2:F57 / G2:G57 columns
# already share their formulas).
$ws.Range("F58:F59").Formula = "=(E58-D58)*24*60"
$ws.Range("G58:G59").Formula = "=F58/60"

# Fix the selection to match the new state.
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("A60").Select()
